$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4041687455978149
$ws.Range("C2").Value = 0.4412235849344484
$ws.Range("D2").Value = 0.07715048283563419
$ws.Range("E2").Value = 0.3965652650683182
$ws.Range("G2").Value = 0.002541939560967374
$ws.Range("I2").Value = 4.411538647597922
$ws.Range("K2").Value = 0.8068701660591273

$ws.Range("B3").Value = 0.3908692952634567
$ws.Range("C3").Value = 0.4042094729968255
$ws.Range("D3").Value = 0.07015876778787344
$ws.Range("E3").Value = 0.345811043711862
$ws.Range("G3").Value = 0.002548262252503955
$ws.Range("I3").Value = 4.070534916726558
$ws.Range("K3").Value = 0.7616356968760556

$ws.Range("B4").Value = 0.3834465772345652
$ws.Range("C4").Value = 0.3816579726803582
$ws.Range("D4").Value = 0.06591116046998025
$ws.Range("E4").Value = 0.3147953914821642
$ws.Range("G4").Value = 0.002552332054598289
$ws.Range("I4").Value = 3.86080993971683
$ws.Range("K4").Value = 0.7347685387720446

$ws.Range("B5").Value = 0.3806063646344171
$ws.Range("C5").Value = 0.3725108200515876
$ws.Range("D5").Value = 0.06419137647291961
$ws.Range("E5").Value = 0.3021901846203576
$ws.Range("G5").Value = 0.002554037928637675
$ws.Range("I5").Value = 3.775244799034169
$ws.Range("K5").Value = 0.724044858176029

$ws.Range("B6").Value = 0.3801458334133088
$ws.Range("C6").Value = 0.3709944873466782
$ws.Range("D6").Value = 0.06390647393328663
$ws.Range("E6").Value = 0.3000990559146004
$ws.Range("G6").Value = 0.00255432405652453
$ws.Range("I6").Value = 3.761030320690196
$ws.Range("K6").Value = 0.7222776915704401

$ws.Range("B7").Value = 0.383407528648263
$ws.Range("C7").Value = 0.3815344395421221
$ws.Range("D7").Value = 0.06588792203864102
$ws.Range("E7").Value = 0.314625260097273
$ws.Range("G7").Value = 0.002552354868485636
$ws.Range("I7").Value = 3.859656402561285
$ws.Range("K7").Value = 0.7346230089670769

$ws.Range("B8").Value = 0.3994274539385003
$ws.Range("C8").Value = 0.4284239213243666
$ws.Range("D8").Value = 0.07473017681671479
$ws.Range("E8").Value = 0.3790326083024098
$ws.Range("G8").Value = 0.002544080807889554
$ws.Range("I8").Value = 4.294024567563184
$ws.Range("K8").Value = 0.791083135073734

$ws.Range("B9").Value = 0.4368466839879659
$ws.Range("C9").Value = 0.5218301252821789
$ws.Range("D9").Value = 0.09244186076475103
$ws.Range("E9").Value = 0.5066588854961509
$ws.Range("G9").Value = 0.00252933452624865
$ws.Range("I9").Value = 5.143701785016901
$ws.Range("K9").Value = 0.9091500910551815

$ws.Range("B10").Value = 0.4681541170267565
$ws.Range("C10").Value = 0.591441561245972
$ws.Range("D10").Value = 0.1057001517268787
$ws.Range("E10").Value = 0.601454866448492
$ws.Range("G10").Value = 0.002519388396983679
$ws.Range("I10").Value = 5.767616627688227
$ws.Range("K10").Value = 1.000598708895353

$ws.Range("B11").Value = 0.4832599991705138
$ws.Range("C11").Value = 0.6233458687815414
$ws.Range("D11").Value = 0.1117892352712744
$ws.Range("E11").Value = 0.6448531390765595
$ws.Range("G11").Value = 0.002515053485096922
$ws.Range("I11").Value = 6.051589242589387
$ws.Range("K11").Value = 1.043273334301631

$ws.Range("B12").Value = 0.4891072995795014
$ws.Range("C12").Value = 0.6354630944701398
$ws.Range("D12").Value = 0.114103653554821
$ws.Range("E12").Value = 0.6613304853457578
$ws.Range("G12").Value = 0.002513439004070989
$ws.Range("I12").Value = 6.159160410733534
$ws.Range("K12").Value = 1.059591630873911

$ws.Range("B13").Value = 0.4878422838467031
$ws.Range("C13").Value = 0.6328518173332327
$ws.Range("D13").Value = 0.1136048141749484
$ws.Range("E13").Value = 0.6577798039327831
$ws.Range("G13").Value = 0.002513785511584606
$ws.Range("I13").Value = 6.135991147566415
$ws.Range("K13").Value = 1.056070088604827

$ws.Range("B14").Value = 0.4837384979584556
$ws.Range("C14").Value = 0.624342034450649
$ws.Range("D14").Value = 0.1119794694619003
$ws.Range("E14").Value = 0.6462078491076255
$ws.Range("G14").Value = 0.002514920119978965
$ws.Range("I14").Value = 6.060438359283637
$ws.Range("K14").Value = 1.044612654401448

$ws.Range("B15").Value = 0.4812414378610299
$ws.Range("C15").Value = 0.6191342507178774
$ws.Range("D15").Value = 0.1109850297966659
$ws.Range("E15").Value = 0.6391254541543248
$ws.Range("G15").Value = 0.002515618617514212
$ws.Range("I15").Value = 6.01416535861307
$ws.Range("K15").Value = 1.037615382731985

$ws.Range("B16").Value = 0.4671845345131658
$ws.Range("C16").Value = 0.5893614567635836
$ws.Range("D16").Value = 0.1053034047558867
$ws.Range("E16").Value = 0.5986245460303508
$ws.Range("G16").Value = 0.002519675491195361
$ws.Range("I16").Value = 5.749062675958299
$ws.Range("K16").Value = 0.9978317660740572

$ws.Range("B17").Value = 0.4587843280731079
$ws.Range("C17").Value = 0.5711588037086699
$ws.Range("D17").Value = 0.1018329319788194
$ws.Range("E17").Value = 0.5738516882330487
$ws.Range("G17").Value = 0.002522212670333718
$ws.Range("I17").Value = 5.586480088260544
$ws.Range("K17").Value = 0.9737036559429555

$ws.Range("B18").Value = 0.4540339014910728
$ws.Range("C18").Value = 0.5607113975324296
$ws.Range("D18").Value = 0.09984223174527074
$ws.Range("E18").Value = 0.5596286568840441
$ws.Range("G18").Value = 0.002523689850345491
$ws.Range("I18").Value = 5.492980130241676
$ws.Range("K18").Value = 0.9599266750415154

$ws.Range("B19").Value = 0.4524393423431547
$ws.Range("C19").Value = 0.5571778643750349
$ws.Range("D19").Value = 0.09916913740539712
$ws.Range("E19").Value = 0.5548172758718124
$ws.Range("G19").Value = 0.002524193072573207
$ws.Range("I19").Value = 5.461324466969756
$ws.Range("K19").Value = 0.9552792341559382

$ws.Range("B20").Value = 0.4596701256788549
$ws.Range("C20").Value = 0.5730941901973097
$ws.Range("D20").Value = 0.1022018062277539
$ws.Range("E20").Value = 0.5764861181438619
$ws.Range("G20").Value = 0.002521940736231144
$ws.Range("I20").Value = 5.603785801946231
$ws.Range("K20").Value = 0.9762616670724071

$ws.Range("B21").Value = 0.4849404091590657
$ws.Range("C21").Value = 0.6268405807916224
$ws.Range("D21").Value = 0.112456636378397
$ws.Range("E21").Value = 0.6496056051438046
$ws.Range("G21").Value = 0.002514586125358205
$ws.Range("I21").Value = 6.082628920024831
$ws.Range("K21").Value = 1.047973651798003

$ws.Range("B22").Value = 0.5021977803743312
$ws.Range("C22").Value = 0.662176131861429
$ws.Range("D22").Value = 0.119209117539171
$ws.Range("E22").Value = 0.6976479975276391
$ws.Range("G22").Value = 0.002509937064952601
$ws.Range("I22").Value = 6.395802980647716
$ws.Range("K22").Value = 1.09576639201785

$ws.Range("B23").Value = 0.4929183820058824
$ws.Range("C23").Value = 0.6432972128753249
$ws.Range("D23").Value = 0.1156004824001116
$ws.Range("E23").Value = 0.6719822935570932
$ws.Range("G23").Value = 0.002512404006063851
$ws.Range("I23").Value = 6.228630618855846
$ws.Range("K23").Value = 1.070172578906408

$ws.Range("B24").Value = 0.4592694108364981
$ws.Range("C24").Value = 0.5722191474461056
$ws.Range("D24").Value = 0.1020350241211503
$ws.Range("E24").Value = 0.5752950328185165
$ws.Range("G24").Value = 0.002522063619792146
$ws.Range("I24").Value = 5.595961980460913
$ws.Range("K24").Value = 0.97510489590303

$ws.Range("B25").Value = 0.4260651907637225
$ws.Range("C25").Value = 0.4963953668455474
$ws.Range("D25").Value = 0.08760864306159988
$ws.Range("E25").Value = 0.47196823464121
$ws.Range("G25").Value = 0.002533166854326638
$ws.Range("I25").Value = 4.913971345221341
$ws.Range("K25").Value = 0.8764002133349038
